# refactor : UUID to AI id
#
# The workbook's "Mulai Semester" column (N) holds the semester code
# "20221" for the sample/template rows. This edit bumps that code to
# "20241" everywhere it appears, and normalizes most of the sample rows
# (N4:N10, N12:N101) from a literal number to the same shared text value
# already used by the two rows (N2, N11) that were stored as text.
# N3 keeps its original numeric cell type, just with the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$col = 14 # column N ("Mulai Semester")

# N2 is already a text cell holding "20221" - this is the template cell
# whose (shared) text every other text cell in the column should match.
$templateCell = $ws.Cells.Item(2, $col)

# Convert the numeric sample cells N4:N10 and N12:N101 into text cells
# carrying the same string as N2/N11, by copying the template cell onto
# each of them (keeps them sharing the exact same text value/type).
for ($r = 4; $r -le 10; $r++) {
    $templateCell.Copy($ws.Cells.Item($r, $col))
}
for ($r = 12; $r -le 101; $r++) {
    $templateCell.Copy($ws.Cells.Item($r, $col))
}

# N3 stays a plain number, just updated to the new semester code.
$ws.Cells.Item(3, $col).Value2 = 20241

# Finally, update the semester code text itself (20221 -> 20241) on every
# cell that carries it, including the two original text cells (N2, N11).
for ($r = 2; $r -le 101; $r++) {
    if ($r -ne 3) {
        $ws.Cells.Item($r, $col).Value2 = "20241"
    }
}

# Restore the view/selection state recorded for the sheet: scrolled so
# column H is the left-most visible column, with N6:N101 selected.
$excel.Goto($ws.Range("H1"), $true)
$null = $ws.Range("N6:N101").Select()
